$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 611-612, pushing the existing 611:626 block
# down to 613:628 (dates/values stay attached to their original rows).
$ws.Range("A611:R612").Insert()

# New row 611 - weekly update, "Primera" quality, Región de Ñuble origin
$ws.Range("A611").Value = 7
$ws.Range("B611").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C611").Value = "Ñuble"
$ws.Range("D611").Value = 45239
$ws.Range("E611").Value = 16
$ws.Range("F611").Value = 100112009
$ws.Range("G611").Value = "Acelga"
$ws.Range("H611").Value = "Sin especificar"
$ws.Range("I611").Value = "Primera"
$ws.Range("J611").Value = 400
$ws.Range("K611").Value = 600
$ws.Range("L611").Value = 700
$ws.Range("M611").Value = 650
$ws.Range("N611").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O611").Value = "Región de Ñuble"
$ws.Range("P611").Value = 650
$ws.Range("Q611").Value = 1
$ws.Range("R611").Value = "Hortaliza"

# New row 612 - weekly update, "Segunda" quality, Región de Ñuble origin
$ws.Range("A612").Value = 7
$ws.Range("B612").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C612").Value = "Ñuble"
$ws.Range("D612").Value = 45239
$ws.Range("E612").Value = 16
$ws.Range("F612").Value = 100112009
$ws.Range("G612").Value = "Acelga"
$ws.Range("H612").Value = "Sin especificar"
$ws.Range("I612").Value = "Segunda"
$ws.Range("J612").Value = 250
$ws.Range("K612").Value = 500
$ws.Range("L612").Value = 500
$ws.Range("M612").Value = 500
$ws.Range("N612").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O612").Value = "Región de Ñuble"
$ws.Range("P612").Value = 500
$ws.Range("Q612").Value = 1
$ws.Range("R612").Value = "Hortaliza"
